$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.826.96"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.58%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.874.40"
$c.NumberFormat = "General"
$c.Style = "Normal"

$ws.Range("E4").Value = "  +0.07%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "301.53"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.01%  "

$ws.Range("E6").Value = "  +0.10%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5355"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.99%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3745"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.90%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07193"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "21.61"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.02%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.8906"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.71%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08201"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.42%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.879.63"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +4.17%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "93.29"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.13%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.314"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.06%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("E17").Value = "  +0.47%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008536"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.70%  "

$ws.Range("E19").Value = "  +0.06%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "26.864.09"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.57%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.992"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.69%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.60"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.03%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "6.391"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.287"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.32%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "146.46"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.37%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "18.09"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.88%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.733"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.13%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "114.09"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.52%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "4.714"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.90%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.618"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.37%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.09112"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.32%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.8088"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.03%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05016"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.25%  "

$ws.Range("E34").Value = "  -4.71%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.958"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.05%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.6118"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +6.30%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.651"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -3.28%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.206"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.77%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01957"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.27%  "

$ws.Range("E40").Value = "  -1.14%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.600"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.28%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "8.874"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.08%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.5151"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +4.35%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "114.89"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.72%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.1496"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.73%  "

$ws.Range("E46").Value = "  +0.12%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.998"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("E48").Value = "  -0.12%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "37.56"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.87%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06076"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.37%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "62.17"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.95%  "
